$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Creación"
$ws.Range("D3").Value = "Publicaciones"
$ws.Range("D4").Value = "Creación"
$ws.Range("D5").Value = "Seminarios"
$ws.Range("D6").Value = "Cursos"

$ws.Range("D7").Select()
